$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.813.20"
$ws.Range("E2").Value = "  +5.23%  "
$ws.Range("D3").Value = "2.273.66"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'232.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").Value = "'0.639"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("D7").Value = "'64.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.29%  "
$ws.Range("E9").Value = "  +7.38%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +16.71%  "
$ws.Range("D11").Value = "'57.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "'26.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.13%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "2.608.47"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "'15.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("D17").Value = "'0.826"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").Value = "2.265.82"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "43.699.61"
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").Value = "0.0₃0996"
$ws.Range("E20").Value = "  +10.82%  "
$ws.Range("D21").Value = "'74.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "'250.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.99%  "
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").Value = "'9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("D28").Value = "'173.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").Value = "'20.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +10.68%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Value = "'0.0686"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("D35").Value = "'5.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "'4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("D37").Value = "'3.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.05%  "
$ws.Range("D38").Value = "'6.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.99%  "
$ws.Range("D39").Value = "'2.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  +5.60%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").Value = "'17.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.61%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0966"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "'10.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.01%  "
$ws.Range("D48").Value = "'98.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "1.478.81"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "'2.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.60%  "
$ws.Range("E51").Value = "  +1.14%  "
